# map layer display revisions; new NWS layer tests
#
# The sheet's "Local Hosted Files" block (rows 22-30) gets a brand-new
# layer entry ("NWS Fire Hazards") inserted immediately after the
# "Fire Weather Forecast" row (row 21), pushing everything from the old
# row 22 down by one row. We replicate that with a native row insert at
# row 22, then populate the freshly-inserted row with the new layer's
# data. Three existing hyperlinks living in column H below the insertion
# point (old H25/H26/H27) need to move down one row to H26/H27/H28 so
# they keep tracking the same description cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at 22, shifting rows 22:30 down to 23:31 -------
$ws.Rows("22:22").Insert()

# --- 2. Fill the new row 22 with the "NWS Fire Hazards" layer entry -----
# (values are set in the same order the original author's sharedStrings
# table shows them being introduced: A, E, H, then B, F, G, then the
# reused C/D strings)
$ws.Range("A22").Value = "NWS Fire Hazards"
$ws.Range("E22").Value = "https://nowcoast.noaa.gov/arcgis/rest/services/nowcoast/wwa_meteoceanhydro_longduration_hazards_time/MapServer"
$ws.Range("H22").Value = "https://nowcoast.noaa.gov/help/#!section=updateschedule"
$ws.Range("B22").Value = "NWS Nowcast disolved polygons - hazardous wildfire conditions"
$ws.Range("F22").Value = "wildfire poly(39)"
$ws.Range("G22").Value = "other products avail; dissolved polygons"
$ws.Range("C22").Value = "NOAA"
$ws.Range("D22").Value = "ArcGIS REST Tile"

# --- 3. Re-anchor the three hyperlinks that sat below the insertion -----
# point. The row-insert operation shifts cell data/styles automatically,
# but not the <hyperlinks> anchors, so those need to be moved by hand:
#   old H25 (SWSI) -> H26, old H26 (CoCoRaHS) -> H27, old H27 (CoCoRaHS) -> H28
$moves = @(
    @{ From = '$H$25'; To = 'H26'; Url = 'https://data.colorado.gov/Water/DWR-Surface-Water-Supply-Index-by-HUC/m9cg-gqek/data' },
    @{ From = '$H$26'; To = 'H27'; Url = 'https://www.cocorahs.org/maps/conditionmonitoring/' },
    @{ From = '$H$27'; To = 'H28'; Url = 'https://www.cocorahs.org/maps/conditionmonitoring/' }
)

foreach ($hl in @($ws.Hyperlinks)) {
    $addr = $hl.Range.Address()
    foreach ($m in $moves) {
        if ($addr -eq $m.From) {
            $hl.Delete()
        }
    }
}

# Remember the pre-existing "hyperlink text" cell style so we can restore
# it after Hyperlinks.Add() (which otherwise mints a fresh, redundant
# style slot instead of reusing the sheet's existing Hyperlink style).
$hyperlinkStyle = $ws.Range("H10").Style
foreach ($m in $moves) {
    $ws.Hyperlinks.Add($ws.Range($m.To), $m.Url, [Type]::Missing, [Type]::Missing, [Type]::Missing) | Out-Null
    $ws.Range($m.To).Style = $hyperlinkStyle
}

# --- 4. Restore the view state (scroll/selection) recorded in the file --
$ws.Range("B12").Select()
